# Update countries & provincias Spain
# Applies the data refresh described by the commit:
#  - A few countries changed ranking position (their row content
#    shifted relative to neighbours) and several numeric values were
#    updated to reflect new case counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos - refreshed totals
Set-Row 4 "Estados Unidos" 1367963 325 256336 1030840 16514 0 80787

# Banglades moves above Ucrania/Rumania with fresh data; Ucrania and
# Rumania shift down one row, keeping their previous values.
Set-Row 37 "Banglades" 15691 1034 2902 12550 1 11 239
Set-Row 38 "Ucrania"   15648 416  3288 11952 207 17 408
Set-Row 39 "Rumania"   15362 0    7051 7350  242 0  961

# Dinamarca - refreshed totals
Set-Row 44 "Dinamarca" 10513 84 8217 1767 40 0 529

# Malasia - refreshed totals
Set-Row 54 "Malasia" 6726 70 5113 1504 20 1 109

# Ghana moves above Afganistan/Nigeria with fresh data; Afganistan and
# Nigeria shift down one row, keeping their previous values.
Set-Row 62 "Ghana"      4700 437 494 4184 5 0 22
Set-Row 63 "Afganistan" 4402 0   558 3724 7 0 120
Set-Row 64 "Nigeria"    4399 0   778 3478 4 0 143

# Vietnam - refreshed totals
Set-Row 136 "Vietnam" 288 0 249 39 2 0 0

# Nueva Caledonia and Belice swap places (rows 192/193)
Set-Row 192 "Nueva Caledonia" 18 0 18 0 0 0 0
Set-Row 193 "Belice"          18 0 16 0 0 0 2

# Curazao and Dominica swap places (rows 198/199)
Set-Row 198 "Curazao"  16 0 14 1 0 0 1
Set-Row 199 "Dominica" 16 0 15 1 0 0 0
